$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the country labels for Albania (row 109) and Mali (row 110):
# the updated source data now lists Mali's (higher) totals in the row
# that used to be Albania's, and Albania keeps its old totals one row
# below - i.e. their labels swap places.
$ws.Range("A109").Value = "Mali"
$ws.Range("A110").Value = "Albania"

# Update the "last updated" timestamp string.
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 21:05"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1656694
$ws.Range("C4").Value = 11600
$ws.Range("E4").Value = 1119878
$ws.Range("G4").Value = 607
$ws.Range("H4").Value = 98254

# Row 5: Brasil
$ws.Range("B5").Value = 340887
$ws.Range("C5").Value = 9997
$ws.Range("E5").Value = 183779

# Row 7: España
$ws.Range("B7").Value = 282370
$ws.Range("C7").Value = 466
$ws.Range("E7").Value = 56734

# Row 10: Francia
$ws.Range("B10").Value = 182469
$ws.Range("C10").Value = 250
$ws.Range("D10").Value = 64547
$ws.Range("E10").Value = 89590
$ws.Range("G10").Value = 43
$ws.Range("H10").Value = 28332

# Row 32: Emiratos Arabes Unidos
$ws.Range("B32").Value = 28704
$ws.Range("C32").Value = 812
$ws.Range("D32").Value = 14495
$ws.Range("E32").Value = 13965
$ws.Range("G32").Value = 3
$ws.Range("H32").Value = 244

# Row 64: Ghana
$ws.Range("D64").Value = 1978
$ws.Range("E64").Value = 4608

# Row 109: now Mali (new, higher totals)
$ws.Range("B109").Value = 1015
$ws.Range("C109").Value = 46
$ws.Range("D109").Value = 574
$ws.Range("E109").Value = 378
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = 63

# Row 110: now Albania (keeps old totals)
$ws.Range("B110").Value = 989
$ws.Range("C110").Value = 8
$ws.Range("D110").Value = 783
$ws.Range("E110").Value = 175
$ws.Range("H110").Value = 31

# Row 116: Costa Rica
$ws.Range("B116").Value = 918
$ws.Range("C116").Value = 7
$ws.Range("D116").Value = 607
